$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 215 (shifts old rows 215..227 down to 216..228)
$ws.Rows.Item(215).Insert()

# Populate the newly inserted row 215 with the new record
$ws.Cells.Item(215, 1).Value = 7
$ws.Cells.Item(215, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(215, 3).Value = "Ñuble"
$ws.Cells.Item(215, 4).Value = 44615
$ws.Cells.Item(215, 5).Value = 16
$ws.Cells.Item(215, 6).Value = 100112002
$ws.Cells.Item(215, 7).Value = "Pimiento"
$ws.Cells.Item(215, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 100
$ws.Cells.Item(215, 11).Value = 6500
$ws.Cells.Item(215, 12).Value = 7000
$ws.Cells.Item(215, 13).Value = 6750
$ws.Cells.Item(215, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(215, 15).Value = "Región del Maule"
$ws.Cells.Item(215, 16).Value = 450
$ws.Cells.Item(215, 17).Value = 15
$ws.Cells.Item(215, 18).Value = "Hortaliza"
